$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 updates
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2: move value from D2 to C2
$ws.Range("D2").ClearContents()
$ws.Range("C2").Value = 28.980596743227572

# Row 3: clear B3 and C3
$ws.Range("B3").ClearContents()
$ws.Range("C3").ClearContents()

# Update selection to reflect B1:E3
$ws.Range("B1:E3").Select()

$wb.Save()
